# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the Leve profit-tracking sheets (ALC, ARM, BSM, CUL, GSM, LTW, WVR)
# per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 742.1707
$ws.Range("I15").Value = 742.1707
$ws.Range("K15").Value = 2226.5121
$ws.Range("M15").Value = -2057.5121
$ws.Range("H40").Value = 7107
$ws.Range("I40").Value = 2874
$ws.Range("J40").Value = 8800.2
$ws.Range("K40").Value = 2874
$ws.Range("L40").Value = 8800.2
$ws.Range("M40").Value = -2699
$ws.Range("N40").Value = -9150.2
$ws.Range("H41").Value = 716.4286
$ws.Range("I41").Value = 213
$ws.Range("J41").Value = 1219.8572
$ws.Range("K41").Value = 213
$ws.Range("L41").Value = 1219.8572
$ws.Range("M41").Value = 227
$ws.Range("N41").Value = -2099.8572
$ws.Range("H58").Value = 976.875
$ws.Range("J58").Value = 1000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3300
$ws.Range("H111").Value = 2172.5386
$ws.Range("I111").Value = 2883.4285
$ws.Range("J111").Value = 1343.1666
$ws.Range("K111").Value = 8650.2855
$ws.Range("L111").Value = 4029.4998
$ws.Range("M111").Value = -5583.2855
$ws.Range("N111").Value = -10163.4998
$ws.Range("H137").Value = 3057.3333
$ws.Range("I137").Value = 1727.8572
$ws.Range("J137").Value = 4220.625
$ws.Range("K137").Value = 5183.571599999999
$ws.Range("L137").Value = 12661.875
$ws.Range("M137").Value = -2633.571599999999
$ws.Range("N137").Value = -17761.875
$ws.Range("H141").Value = 799
$ws.Range("I141").Value = 799
$ws.Range("K141").Value = 2397
$ws.Range("M141").Value = 2783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1474.75
$ws.Range("I63").Value = 1499.6666
$ws.Range("K63").Value = 1499.6666
$ws.Range("M63").Value = -813.6666
$ws.Range("H66").Value = 1474.75
$ws.Range("I66").Value = 1499.6666
$ws.Range("K66").Value = 7498.333000000001
$ws.Range("M66").Value = -4066.333000000001
$ws.Range("H74").Value = 2829.5557
$ws.Range("I74").Value = 2829.5557
$ws.Range("K74").Value = 2829.5557
$ws.Range("M74").Value = -1955.5557
$ws.Range("H77").Value = 2829.5557
$ws.Range("I77").Value = 2829.5557
$ws.Range("K77").Value = 14147.7785
$ws.Range("M77").Value = -9779.7785
$ws.Range("H88").Value = 6141.5713
$ws.Range("I88").Value = 5745.75
$ws.Range("K88").Value = 5745.75
$ws.Range("M88").Value = -5339.75
$ws.Range("H91").Value = 6141.5713
$ws.Range("I91").Value = 5745.75
$ws.Range("K91").Value = 5745.75
$ws.Range("M91").Value = -4341.75
$ws.Range("H132").Value = 1021.4286
$ws.Range("I132").Value = 1021.4286
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3064.2858
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -534.2857999999997
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4659.2666
$ws.Range("I105").Value = 3686.25
$ws.Range("K105").Value = 3686.25
$ws.Range("M105").Value = -1939.25
$ws.Range("H134").Value = 2766.3333
$ws.Range("I134").Value = 2766.3333
$ws.Range("K134").Value = 8298.999899999999
$ws.Range("M134").Value = -5763.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 182.7
$ws.Range("I11").Value = 195.66667
$ws.Range("J11").Value = 177.14285
$ws.Range("K11").Value = 587.00001
$ws.Range("L11").Value = 531.4285500000001
$ws.Range("M11").Value = -447.00001
$ws.Range("N11").Value = -811.4285500000001
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330
$ws.Range("H122").Value = 1496.375
$ws.Range("J122").Value = 1443.3636
$ws.Range("L122").Value = 12990.2724
$ws.Range("N122").Value = -17890.2724
$ws.Range("H128").Value = 499996.5
$ws.Range("I128").Value = 499996.5
$ws.Range("K128").Value = 1499989.5
$ws.Range("M128").Value = -1495009.5
$ws.Range("H137").Value = 1030
$ws.Range("I137").Value = 1030
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3090
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 2010
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 3496.3333
$ws.Range("I139").Value = 1779.8572
$ws.Range("J139").Value = 4998.25
$ws.Range("K139").Value = 5339.571599999999
$ws.Range("L139").Value = 14994.75
$ws.Range("M139").Value = -199.5715999999993
$ws.Range("N139").Value = -25274.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33336004
$ws.Range("I70").Value = 33336004
$ws.Range("K70").Value = 33336004
$ws.Range("M70").Value = -33335734
$ws.Range("H73").Value = 33336004
$ws.Range("I73").Value = 33336004
$ws.Range("K73").Value = 33336004
$ws.Range("M73").Value = -33335068
$ws.Range("H102").Value = 944.8
$ws.Range("I102").Value = 944.8
$ws.Range("K102").Value = 944.8
$ws.Range("M102").Value = 677.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H16").Value = 1070.85
$ws.Range("I16").Value = 1169.3846
$ws.Range("J16").Value = 887.8571
$ws.Range("K16").Value = 1169.3846
$ws.Range("L16").Value = 887.8571
$ws.Range("M16").Value = -999.3846000000001
$ws.Range("N16").Value = -1227.8571
$ws.Range("H40").Value = 14257.818
$ws.Range("I40").Value = 14104.5
$ws.Range("K40").Value = 14104.5
$ws.Range("M40").Value = -13968.5
$ws.Range("H46").Value = 1305.4445
$ws.Range("I46").Value = 1016.3333
$ws.Range("J46").Value = 1450
$ws.Range("K46").Value = 1016.3333
$ws.Range("L46").Value = 1450
$ws.Range("M46").Value = -828.3333
$ws.Range("N46").Value = -1826
$ws.Range("H54").Value = 100000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 100000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 100000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -101288
$ws.Range("H61").Value = 3999.5
$ws.Range("I61").Value = 3999.5
$ws.Range("K61").Value = 3999.5
$ws.Range("M61").Value = -3797.5
$ws.Range("H113").Value = 3999.5
$ws.Range("I113").Value = 3999.5
$ws.Range("K113").Value = 3999.5
$ws.Range("M113").Value = -1829.5
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 4749.25
$ws.Range("I132").Value = 2665.6667
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 7997.000100000001
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = -5467.000100000001
$ws.Range("N132").Value = -38060
$ws.Range("H136").Value = 4750
$ws.Range("J136").Value = 5250
$ws.Range("L136").Value = 15750
$ws.Range("N136").Value = -20850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 826.875
$ws.Range("I113").Value = 558
$ws.Range("J113").Value = 1275
$ws.Range("K113").Value = 1674
$ws.Range("L113").Value = 3825
$ws.Range("M113").Value = 496
$ws.Range("N113").Value = -8165
$ws.Range("H126").Value = 1789.4286
$ws.Range("I126").Value = 1486.6364
$ws.Range("J126").Value = 2899.6667
$ws.Range("K126").Value = 4459.9092
$ws.Range("L126").Value = 8699.000100000001
$ws.Range("M126").Value = -1989.9092
$ws.Range("N126").Value = -13639.0001
$ws.Range("H132").Value = 4305
$ws.Range("I132").Value = 4305
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12915
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10385
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 24165.39
$ws.Range("I136").Value = 24809.273
$ws.Range("K136").Value = 74427.819
$ws.Range("M136").Value = -71877.819

